$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '64.348.65'
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '  +0.03%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '3.120.15'
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '  -5.61%  '
$ws.Range('E4').Value = '  -0.75%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '591.25'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +0.58%  '
$ws.Range('E6').Value = '  +4.42%  '
$ws.Range('E7').Value = '  -0.93%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '3.117.69'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  -0.99%  '
$ws.Range('E9').Value = '  +0.47%  '
$ws.Range('E10').Value = '  -0.79%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '5.94'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  +1.83%  '
$ws.Range('E12').Value = '  +0.77%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '37.98'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  +2.24%  '
$ws.Range('E14').Value = '  -1.23%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '3.636.44'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  -5.45%  '
$ws.Range('E16').Value = '  -1.62%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '7.24'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  +2.24%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '64.050.38'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  -0.10%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '3.116.98'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  -3.26%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '470.71'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +1.19%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '14.92'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +4.32%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '0.739'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +0.97%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '7.63'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +2.58%  '
$ws.Range('B24').Value = 'InternetComputer(DFINITY)'
$ws.Range('C24').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '13.33'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +2.99%  '
$ws.Range('B25').Value = 'Fetch.AI'
$ws.Range('C25').Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '2.39'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +6.57%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '81.89'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  +0.92%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '0.998'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -0.58%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '9.93'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +6.89%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '7.46'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  +4.85%  '
$ws.Range('E30').Value = '  +1.14%  '
$ws.Range('E31').Value = '  -0.73%  '
$ws.Range('E32').Value = '  -0.34%  '
$ws.Range('E33').Value = '  +6.24%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '27.62'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  +2.07%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.0₃0857'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  -0.46%  '
$ws.Range('E36').Value = '  +1.46%  '
$ws.Range('E37').Value = '  +2.82%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '6.17'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +2.44%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '2.27'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  -2.39%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '9.41'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +6.52%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '456.78'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +4.05%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '51.07'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -0.70%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.290'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  +0.66%  '
$ws.Range('E44').Value = '  -0.02%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '2.861.70'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -1.98%  '
$ws.Range('E46').Value = '  +1.96%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '38.81'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -1.64%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '129.73'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +2.50%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '25.42'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  +5.16%  '
$ws.Range('E50').Value = '  +4.39%  '
$ws.Range('E51').Value = '  +0.03%  '
